# Insert a new LOINC concept row ("34111-5") into the "Include from LOINC"
# sheet, just above the existing "System URI" row, pushing that row (and the
# blank row above it) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Include from LOINC")

# Row 11 currently holds a blank spacer row ("", "") immediately followed by
# the "System URI" / "http://loinc.org" row at 12. Insert a fresh row at 11
# so the new concept code lands there and everything below shifts down one.
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = "34111-5"
$ws.Cells.Item(11, 2).Value = ""

# Row insertion doesn't inherit the bordered/top-aligned data-row style used
# by the rest of the table, so copy that formatting down from the row above.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
